# Quarterly indexing bug-fix: shift column A date serials for rows 2-63
# from the first day of a calendar quarter (Jan/Apr/Jul/Oct 1) to the
# 15th day of the middle month of that quarter (Feb/May/Aug/Nov 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 25614
$ws.Range("A3").Value = 25614
$ws.Range("A4").Value = 25614
$ws.Range("A5").Value = 25614
$ws.Range("A6").Value = 25614
$ws.Range("A7").Value = 25614
$ws.Range("A8").Value = 25614
$ws.Range("A9").Value = 25614
$ws.Range("A10").Value = 25614
$ws.Range("A11").Value = 25614
$ws.Range("A12").Value = 40313
$ws.Range("A13").Value = 40405
$ws.Range("A14").Value = 40497
$ws.Range("A15").Value = 40589
$ws.Range("A16").Value = 40678
$ws.Range("A17").Value = 40770
$ws.Range("A18").Value = 40862
$ws.Range("A19").Value = 40954
$ws.Range("A20").Value = 41044
$ws.Range("A21").Value = 41136
$ws.Range("A22").Value = 41228
$ws.Range("A23").Value = 41320
$ws.Range("A24").Value = 41409
$ws.Range("A25").Value = 41501
$ws.Range("A26").Value = 41593
$ws.Range("A27").Value = 41685
$ws.Range("A28").Value = 41774
$ws.Range("A29").Value = 41866
$ws.Range("A30").Value = 41958
$ws.Range("A31").Value = 42050
$ws.Range("A32").Value = 42139
$ws.Range("A33").Value = 42231
$ws.Range("A34").Value = 42323
$ws.Range("A35").Value = 42415
$ws.Range("A36").Value = 42505
$ws.Range("A37").Value = 42597
$ws.Range("A38").Value = 42689
$ws.Range("A39").Value = 42781
$ws.Range("A40").Value = 42870
$ws.Range("A41").Value = 42962
$ws.Range("A42").Value = 43054
$ws.Range("A43").Value = 43146
$ws.Range("A44").Value = 43235
$ws.Range("A45").Value = 43327
$ws.Range("A46").Value = 43419
$ws.Range("A47").Value = 43511
$ws.Range("A48").Value = 43600
$ws.Range("A49").Value = 43692
$ws.Range("A50").Value = 43784
$ws.Range("A51").Value = 43876
$ws.Range("A52").Value = 43966
$ws.Range("A53").Value = 44058
$ws.Range("A54").Value = 44150
$ws.Range("A55").Value = 44242
$ws.Range("A56").Value = 44331
$ws.Range("A57").Value = 44423
$ws.Range("A58").Value = 44515
$ws.Range("A59").Value = 44607
$ws.Range("A60").Value = 44696
$ws.Range("A61").Value = 44788
$ws.Range("A62").Value = 44880
$ws.Range("A63").Value = 44972
